$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

# Row 18: Baixa Grande x Lagoa Verde -> 4x6, Finalizado
$ws.Range("E18").Value = "4x6"
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 6
$ws.Range("I18").Value = "Finalizado"

# Row 19: Pró Limp x Es. Po. Seguro -> 2x3, Finalizado
$ws.Range("E19").Value = "2x3"
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = "Finalizado"

# Row 20: At. Ma. Martins x Real Cutias -> 5x2, Finalizado
$ws.Range("E20").Value = "5x2"
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 2
$ws.Range("I20").Value = "Finalizado"

# Move the active selection to L19 (matches the saved view state in the workbook)
$ws.Range("L19").Select() | Out-Null
